# Updates cryptocurrency price/volume/hour data to the latest scrape,
# matching the "Updated symbol list" GitHub Actions commit.
# All affected cells (Price, Volume(1h), Hora columns) are plain text
# values in the original workbook, so we force each cell to the "@"
# (text) number format before writing the new value, then restore the
# cell style to "Normal" so no lingering custom number format is left
# applied to the cell (keeps formatting identical to the original).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue 'D2' '261.47'
Set-TextValue 'E2' '0.07%'
Set-TextValue 'G2' '10'
# Row 3
Set-TextValue 'E3' '-2.47%'
Set-TextValue 'G3' '10'
# Row 4
Set-TextValue 'D4' '4.706'
Set-TextValue 'E4' '-0.26%'
Set-TextValue 'G4' '10'
# Row 5
Set-TextValue 'E5' '1.81%'
Set-TextValue 'G5' '10'
# Row 6
Set-TextValue 'D6' '6.739'
Set-TextValue 'E6' '0.87%'
Set-TextValue 'G6' '10'
# Row 7
Set-TextValue 'D7' '0.8496'
Set-TextValue 'E7' '0.12%'
Set-TextValue 'G7' '10'
# Row 8
Set-TextValue 'D8' '0.9114'
Set-TextValue 'E8' '-1.61%'
Set-TextValue 'G8' '10'
# Row 9
Set-TextValue 'D9' '0.1404'
Set-TextValue 'E9' '0.02%'
Set-TextValue 'G9' '10'
# Row 10
Set-TextValue 'D10' '0.05130'
Set-TextValue 'E10' '9.60%'
Set-TextValue 'G10' '10'
# Row 11
Set-TextValue 'D11' '0.07099'
Set-TextValue 'E11' '-0.04%'
Set-TextValue 'G11' '10'
# Row 12
Set-TextValue 'D12' '0.03108'
Set-TextValue 'E12' '0.74%'
Set-TextValue 'G12' '10'
# Row 13
Set-TextValue 'D13' '0.09050'
Set-TextValue 'E13' '-0.16%'
Set-TextValue 'G13' '10'
# Row 14
Set-TextValue 'D14' '0.001531'
Set-TextValue 'E14' '0.05%'
Set-TextValue 'G14' '10'
# Row 15
Set-TextValue 'D15' '0.0006168'
Set-TextValue 'E15' '1.52%'
Set-TextValue 'G15' '10'
# Row 16
Set-TextValue 'D16' '0.005940'
Set-TextValue 'E16' '-2.15%'
Set-TextValue 'G16' '10'
# Row 17
Set-TextValue 'D17' '3.447'
Set-TextValue 'E17' '-0.04%'
Set-TextValue 'G17' '10'
# Row 18
Set-TextValue 'D18' '3.172'
Set-TextValue 'E18' '0.85%'
Set-TextValue 'G18' '10'
# Row 19
Set-TextValue 'D19' '2.188'
Set-TextValue 'E19' '1.15%'
Set-TextValue 'G19' '10'
# Row 20
Set-TextValue 'D20' '0.3103'
Set-TextValue 'E20' '-0.17%'
Set-TextValue 'G20' '10'
# Row 21
Set-TextValue 'E21' '0.39%'
Set-TextValue 'G21' '10'
# Row 22
Set-TextValue 'D22' '4.095'
Set-TextValue 'E22' '0.19%'
Set-TextValue 'G22' '10'
# Row 23
Set-TextValue 'D23' '0.04229'
Set-TextValue 'E23' '-0.37%'
Set-TextValue 'G23' '10'
# Row 24
Set-TextValue 'D24' '0.001182'
Set-TextValue 'E24' '-3.35%'
Set-TextValue 'G24' '10'
# Row 25
Set-TextValue 'D25' '0.004060'
Set-TextValue 'E25' '7.23%'
Set-TextValue 'G25' '10'
# Row 26
Set-TextValue 'E26' '0.04%'
Set-TextValue 'G26' '10'
# Row 27
Set-TextValue 'E27' '4.10%'
Set-TextValue 'G27' '10'
# Row 28
Set-TextValue 'G28' '10'
# Row 29
Set-TextValue 'G29' '10'
# Row 30
Set-TextValue 'G30' '10'
# Row 31
Set-TextValue 'G31' '10'
# Row 32
Set-TextValue 'G32' '10'
# Row 33
Set-TextValue 'G33' '10'
# Row 34
Set-TextValue 'G34' '10'
# Row 35
Set-TextValue 'G35' '10'
# Row 36
Set-TextValue 'G36' '10'
# Row 37
Set-TextValue 'G37' '10'
# Row 38
Set-TextValue 'G38' '10'
# Row 39
Set-TextValue 'G39' '10'
# Row 40
Set-TextValue 'D40' '0.03943'
Set-TextValue 'E40' '1.81%'
Set-TextValue 'G40' '10'
# Row 41
Set-TextValue 'D41' '0.1113'
Set-TextValue 'E41' '-0.06%'
Set-TextValue 'G41' '10'
# Row 42
Set-TextValue 'D42' '0.004132'
Set-TextValue 'E42' '0.93%'
Set-TextValue 'G42' '10'
# Row 43
Set-TextValue 'D43' '0.002142'
Set-TextValue 'E43' '-3.36%'
Set-TextValue 'G43' '10'
# Row 44
Set-TextValue 'D44' '0.01324'
Set-TextValue 'E44' '-18.61%'
Set-TextValue 'G44' '10'
# Row 45
Set-TextValue 'D45' '0.00005164'
Set-TextValue 'E45' '0.51%'
Set-TextValue 'G45' '10'
# Row 46
Set-TextValue 'E46' '0.04%'
Set-TextValue 'G46' '10'
# Row 47
Set-TextValue 'G47' '10'
# Row 48
Set-TextValue 'D48' '0.2496'
Set-TextValue 'E48' '84.45%'
Set-TextValue 'G48' '10'
# Row 49
Set-TextValue 'E49' '0.04%'
Set-TextValue 'G49' '10'
# Row 50
Set-TextValue 'E50' '0.04%'
Set-TextValue 'G50' '10'
# Row 51
Set-TextValue 'G51' '10'
